$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulation results data (simOut) with simplified/refined values
$ws.Range("B2").Value = 34.827570920459564
$ws.Range("C2").Value = 16.095624835736984
$ws.Range("D2").Value = 0.46215180704094294
$ws.Range("E2").Value = 30.621983369473156
$ws.Range("F2").Value = 15.265674596829825
$ws.Range("G2").Value = 0.49852011258186729
$ws.Range("H2").Value = 338.5
$ws.Range("I2").Value = 302.5

$ws.Range("B3").Value = 35.042927225540694
$ws.Range("C3").Value = 16.288836562883279
$ws.Range("D3").Value = 0.46482522587357716
$ws.Range("E3").Value = 30.466776504714311
$ws.Range("F3").Value = 15.146040905679763
$ws.Range("G3").Value = 0.49713302959163153
$ws.Range("H3").Value = 338.5
$ws.Range("I3").Value = 302

$ws.Range("B4").Value = 34.733189821720366
$ws.Range("C4").Value = 16.040995910440088
$ws.Range("D4").Value = 0.46183480390876358
$ws.Range("E4").Value = 30.601932513401579
$ws.Range("F4").Value = 15.25660732043486
$ws.Range("G4").Value = 0.498550453104669
$ws.Range("H4").Value = 338
$ws.Range("I4").Value = 302

$ws.Range("B5").Value = 34.878446303742066
$ws.Range("C5").Value = 16.114132308030729
$ws.Range("D5").Value = 0.4620083179078382
$ws.Range("E5").Value = 30.729153510854232
$ws.Range("F5").Value = 15.368503922073144
$ws.Range("G5").Value = 0.50012779937607588
$ws.Range("H5").Value = 338
$ws.Range("I5").Value = 301.5

$ws.Range("B6").Value = 35.007017093821759
$ws.Range("C6").Value = 16.233214148377709
$ws.Range("D6").Value = 0.46371314942005271
$ws.Range("E6").Value = 30.489932126339767
$ws.Range("F6").Value = 15.167643760899209
$ws.Range("G6").Value = 0.49746400543135755
$ws.Range("H6").Value = 338
$ws.Range("I6").Value = 301

$ws.Range("B7").Value = 35.100392017634896
$ws.Range("C7").Value = 16.286205171736995
$ws.Range("D7").Value = 0.46398926722968203
$ws.Range("E7").Value = 30.567216220576341
$ws.Range("F7").Value = 15.224066924400196
$ws.Range("G7").Value = 0.49805212272330202
$ws.Range("I7").Value = 301

$ws.Range("B8").Value = 34.691568300702841
$ws.Range("C8").Value = 15.986264184264474
$ws.Range("D8").Value = 0.46081122783776246
$ws.Range("E8").Value = 30.649578796203976
$ws.Range("F8").Value = 15.296118262943807
$ws.Range("G8").Value = 0.49906455043480952
$ws.Range("H8").Value = 338
$ws.Range("I8").Value = 301

$ws.Range("B9").Value = 34.775742751376669
$ws.Range("C9").Value = 16.062318188654977
$ws.Range("D9").Value = 0.46188282169815387
$ws.Range("E9").Value = 30.748877890738044
$ws.Range("F9").Value = 15.397482198186733
$ws.Range("G9").Value = 0.50074940142204838
$ws.Range("H9").Value = 338
$ws.Range("I9").Value = 301

$ws.Range("B10").Value = 34.819224916615333
$ws.Range("C10").Value = 16.101283467201643
$ws.Range("D10").Value = 0.46242509722030878
$ws.Range("E10").Value = 30.806176384232835
$ws.Range("F10").Value = 15.433478317216988
$ws.Range("G10").Value = 0.50098649454971389
$ws.Range("H10").Value = 338
$ws.Range("I10").Value = 300.5

$ws.Range("B11").Value = 34.879840706953978
$ws.Range("C11").Value = 16.115517118376538
$ws.Range("D11").Value = 0.46202955035753918
$ws.Range("E11").Value = 30.860874543032381
$ws.Range("F11").Value = 15.505205216596877
$ws.Range("G11").Value = 0.50242274226468953
$ws.Range("H11").Value = 338
$ws.Range("I11").Value = 300.5

# Adjust column widths
$ws.Columns.Item(1).ColumnWidth = 29.166666666666668
$ws.Columns.Item(2).ColumnWidth = 29.5
$ws.Columns.Item(3).ColumnWidth = 28.0
$ws.Columns.Item(4).ColumnWidth = 34.333333333333336
$ws.Columns.Item(5).ColumnWidth = 28.166666666666668
$ws.Columns.Item(6).ColumnWidth = 26.833333333333332
$ws.Columns.Item(7).ColumnWidth = 33.166666666666664
$ws.Columns.Item(8).ColumnWidth = 28.333333333333332
$ws.Columns.Item(9).ColumnWidth = 27.166666666666668
